$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record needs to be inserted before the existing row 217,
# shifting all subsequent rows (217..320) down by one (to 218..321).
$ws.Rows.Item(217).Insert()

# Match the date number format used throughout column D before writing the date value.
$ws.Range("D217").NumberFormat = $ws.Range("D218").NumberFormat

# Populate the newly inserted row with the new record's data.
$ws.Range("A217").Value = 5
$ws.Range("B217").Value = "Macroferia Regional de Talca"
$ws.Range("C217").Value = "Maule"
$ws.Range("D217").Value = 44523
$ws.Range("E217").Value = 7
$ws.Range("F217").Value = "Fruta"
$ws.Range("G217").Value = 100103
$ws.Range("H217").Value = "Frutos de hueso (carozo)"
$ws.Range("I217").Value = 100103006
$ws.Range("J217").Value = "Nectarín"
$ws.Range("K217").Value = "Early John"
$ws.Range("L217").Value = "Primera"
$ws.Range("M217").Value = 60
$ws.Range("N217").Value = 18000
$ws.Range("O217").Value = 18000
$ws.Range("P217").Value = 18000
$ws.Range("Q217").Value = "$/caja 15 kilos empedrada"
$ws.Range("R217").Value = "Región de O'Higgins"
$ws.Range("S217").Value = 1200
$ws.Range("T217").Value = 15
